$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (shifts MethodName/Desc/Resource/ExtraEntry* one column right)
$ws.Columns("D").Insert() | Out-Null

# New column header: InvokeType
$ws.Range("D1").Value = "InvokeType"

# Text polish on the skill descriptions (now living in column F after the insert)
$ws.Range("F3").Value = "升级你的混元剑坯，御利分别有着不同的效果，二者随机出现其一。\n混元·御：使用混元剑坯时，增加自身15%攻击伤害的护盾，持续一回合\n混元·利：增加混元剑坯10%攻击伤害"
$ws.Range("F4").Value = "增加8护盾"
$ws.Range("F5").Value = "增加混元剑坯5攻击伤害"

# Cosmetic: update the active selection like in the authored workbook
$ws.Range("F14").Select() | Out-Null
